# Refresh the cryptos list (prices + 1h volume deltas) for the Wed Jun 28
# 2023 GitHub Actions run. A handful of rows were also re-ranked, so their
# Coin / Link / Price / Volume columns swap with the neighboring row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Cells D2:D51 store prices as plain text (e.g. "30.326.57" uses two
    # dots and can never be a number, but "1.001" is also text in this
    # sheet and must NOT be silently promoted to the number 1.001 by Excel's
    # auto-detection). Prefix ambiguous-looking numerics with a text
    # quote-prefix, then strip the resulting cell style back to Normal so
    # no stray number-format/quote-prefix style is left behind.
    $range = $ws.Range($cellRef)
    $looksNumeric = $text -match "^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?%?$"
    if ($looksNumeric) {
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

Set-TextCell "D2" "30.167.49"
Set-TextCell "E2" "  -1.75%  "

Set-TextCell "D3" "1.856.87"
Set-TextCell "E3" "  -1.42%  "

Set-TextCell "D4" "1.001"
Set-TextCell "E4" "  +0.14%  "

Set-TextCell "D5" "232.90"
Set-TextCell "E5" "  -2.87%  "

Set-TextCell "D6" "1.002"
Set-TextCell "E6" "  +0.13%  "

Set-TextCell "D7" "0.4733"
Set-TextCell "E7" "  -1.92%  "

Set-TextCell "D8" "0.2731"
Set-TextCell "E8" "  -3.99%  "

Set-TextCell "D9" "0.06408"
Set-TextCell "E9" "  -2.20%  "

Set-TextCell "D10" "1.848.56"
Set-TextCell "E10" "  -4.79%  "

Set-TextCell "D11" "0.07451"
Set-TextCell "E11" "  -0.67%  "

Set-TextCell "D12" "16.20"
Set-TextCell "E12" "  -2.92%  "

Set-TextCell "D13" "4.988"
Set-TextCell "E13" "  -2.77%  "

Set-TextCell "D14" "85.08"
Set-TextCell "E14" "  -4.63%  "

Set-TextCell "D15" "0.6326"
Set-TextCell "E15" "  -5.41%  "

Set-TextCell "D16" "30.102.49"
Set-TextCell "E16" "  -1.85%  "

Set-TextCell "E17" "  +0.15%  "

Set-TextCell "B18" "Avalanche"
Set-TextCell "C18" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell "D18" "12.76"
Set-TextCell "E18" "  -4.84%  "

Set-TextCell "B19" "BitcoinCash"
Set-TextCell "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell "D19" "229.05"
Set-TextCell "E19" "  -1.86%  "

Set-TextCell "D20" "0.000007324"
Set-TextCell "E20" "  -4.20%  "

Set-TextCell "D21" "2.100.94"
Set-TextCell "E21" "  -3.57%  "

Set-TextCell "D22" "1.002"
Set-TextCell "E22" "  +0.18%  "

Set-TextCell "D23" "5.076"
Set-TextCell "E23" "  -4.79%  "

Set-TextCell "D24" "5.986"
Set-TextCell "E24" "  -3.95%  "

Set-TextCell "B25" "Cosmos"
Set-TextCell "C25" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D25" "9.228"
Set-TextCell "E25" "  -1.47%  "

Set-TextCell "B26" "Monero"
Set-TextCell "C26" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D26" "166.60"
Set-TextCell "E26" "  -1.20%  "

Set-TextCell "D27" "17.77"
Set-TextCell "E27" "  -4.76%  "

Set-TextCell "D28" "1.879"
Set-TextCell "E28" "  -4.09%  "

Set-TextCell "D29" "1.384"
Set-TextCell "E29" "  -3.82%  "

Set-TextCell "D30" "0.09989"
Set-TextCell "E30" "  +4.12%  "

Set-TextCell "D31" "4.144"
Set-TextCell "E31" "  -5.06%  "

Set-TextCell "D32" "3.908"
Set-TextCell "E32" "  -3.72%  "

Set-TextCell "D33" "0.04880"
Set-TextCell "E33" "  -3.49%  "

Set-TextCell "D34" "1.138"
Set-TextCell "E34" "  -7.10%  "

Set-TextCell "D35" "0.7169"
Set-TextCell "E35" "  -4.82%  "

Set-TextCell "D36" "1.001"
Set-TextCell "E36" "  -0.38%  "

Set-TextCell "D37" "2.699"
Set-TextCell "E37" "  -0.25%  "

Set-TextCell "D38" "0.01897"
Set-TextCell "E38" "  +1.50%  "

Set-TextCell "D39" "2.632"
Set-TextCell "E39" "  +0.16%  "

Set-TextCell "D40" "0.9003"
Set-TextCell "E40" "  -2.03%  "

Set-TextCell "D41" "1.954"
Set-TextCell "E41" "  -7.36%  "

Set-TextCell "D42" "105.81"
Set-TextCell "E42" "  -0.63%  "

Set-TextCell "D43" "1.000"
Set-TextCell "E43" "  -0.07%  "

Set-TextCell "B44" "FraxShare"
Set-TextCell "C44" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D44" "5.558"
Set-TextCell "E44" "  -4.79%  "

Set-TextCell "B45" "TheSandbox"
Set-TextCell "C45" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell "D45" "0.4085"
Set-TextCell "E45" "  -5.21%  "

Set-TextCell "D46" "7.024"
Set-TextCell "E46" "  -6.12%  "

Set-TextCell "D47" "60.89"
Set-TextCell "E47" "  -6.58%  "

Set-TextCell "D48" "0.1194"
Set-TextCell "E48" "  -7.95%  "

Set-TextCell "D49" "8.686"
Set-TextCell "E49" "  -3.04%  "

Set-TextCell "B50" "Elrond"
Set-TextCell "C50" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextCell "D50" "33.11"
Set-TextCell "E50" "  -2.80%  "

Set-TextCell "B51" "NEARProtocol"
Set-TextCell "C51" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D51" "1.394"
Set-TextCell "E51" "  -6.39%  "
